{"js": "// Problem 3 solution evaluation\n// Insert a blank paragraph, an \"Evaluate each potential solution:\" heading\n// paragraph, and the evaluation text paragraph right after the \"Identify\n// potential solutions\" answer paragraph for Problem 3 (the one that starts\n// with \"First solution: create a massive chart\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its exact text.\nconst anchorText =\n  \"First solution: create a massive chart the would display a column of \" +\n  \"each possible finger and then a row showing what numbers are counted \" +\n  \"for that finger.  Second solution:  create a mathematical formula \" +\n  \"that will help figure out which finger the count will stop on.\";\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the Problem 3 'Identify potential solutions' paragraph.\");\n}\n\n// Insert in reverse order, each time right \"After\" the anchor paragraph,\n// so the final order reading top-to-bottom is:\n//   (anchor) -> \"\" -> \"Evaluate each potential solution:\" -> \"Each solution...\"\nanchor.insertParagraph(\n  \"Each solution would meet the goals.  Each solution will both work in any and all cases.\",\n  \"After\"\n);\nanchor.insertParagraph(\"Evaluate each potential solution:\", \"After\");\nanchor.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Problem 3 solution evaluation\n# Insert a blank paragraph, an \"Evaluate each potential solution:\" heading\n# paragraph, and the evaluation text paragraph right after the \"Identify\n# potential solutions\" answer paragraph for Problem 3 (the one that starts\n# with \"First solution: create a massive chart\").\n\n$d = $word.ActiveDocument\n\n$anchorText = \"First solution: create a massive chart the would display a column of each possible finger and then a row showing what numbers are counted for that finger.  Second solution:  create a mathematical formula that will help figure out which finger the count will stop on.\"\n\n$target = $null\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd() -eq $anchorText) {\n        $target = $p\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the Problem 3 'Identify potential solutions' paragraph.\"\n}\n\n# Create three new empty paragraphs right after the anchor paragraph.\n$r = $target.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n$r.InsertParagraphAfter()\n$r.InsertParagraphAfter()\n\n# Fill in the text for the 2nd and 3rd new paragraphs; the 1st stays blank.\n$d.Paragraphs.Item($targetIndex + 2).Range.Text = \"Evaluate each potential solution:\"\n$d.Paragraphs.Item($targetIndex + 3).Range.Text = \"Each solution would meet the goals.  Each solution will both work in any and all cases.\"\n"}
